$wb = $excel.ActiveWorkbook

# --- Overview sheet: update Status cells for zh-cn (E2) and de-de (F2) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# --- zh-cn sheet: update Status cell (C2) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "In Translation"

# --- de-de sheet: update Status cell (C2) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "In Translation"

# --- Shrink the now-shorter "Status" columns to match the regenerated report's autofit widths ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZh.Columns.Item(3).ColumnWidth = 12.5
$wsDe.Columns.Item(3).ColumnWidth = 12.5
